$wb = $excel.ActiveWorkbook

# ---- Sheet: LP1912 ----
$ws = $wb.Worksheets.Item('LP1912')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 13:18:32'
$ws.Cells.Item(3, 1).Value = 'Total filas: 248'
$ws.Cells.Item(39, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(40, 3).Value = '86_EST CHICA-ESC AGRARIA'
$ws.Cells.Item(63, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(64, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(65, 1).Value = '05:52:07'
$ws.Cells.Item(65, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(65, 4).Value = 100
$ws.Cells.Item(66, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(67, 1).Value = '07:28:14'
$ws.Cells.Item(67, 3).Value = '84_COLONIA URQUIZA-ESC 49'
$ws.Cells.Item(67, 4).Value = 4
$ws.Cells.Item(72, 1).Value = '06:59:37'
$ws.Cells.Item(72, 3).Value = '14_ABASTO'
$ws.Cells.Item(72, 4).Value = 48
$ws.Cells.Item(73, 1).Value = '07:28:14'
$ws.Cells.Item(73, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(73, 4).Value = 19
$ws.Cells.Item(88, 1).Value = '07:28:14'
$ws.Cells.Item(88, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(88, 4).Value = 55
$ws.Cells.Item(89, 1).Value = '08:13:38'
$ws.Cells.Item(89, 3).Value = '215B_EL PATO'
$ws.Cells.Item(89, 4).Value = 10
$ws.Cells.Item(153, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(154, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(163, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(164, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(184, 1).Value = '12:02:21'
$ws.Cells.Item(184, 3).Value = '17_ROMERO'
$ws.Cells.Item(184, 4).Value = 12
$ws.Cells.Item(185, 1).Value = '10:57:35'
$ws.Cells.Item(185, 3).Value = '10_OLMOS'
$ws.Cells.Item(185, 4).Value = 77
$ws.Cells.Item(189, 3).Value = '215A_EL PATO'
$ws.Cells.Item(190, 3).Value = '14_ABASTO'
$ws.Cells.Item(193, 1).Value = '11:49:23'
$ws.Cells.Item(193, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(193, 4).Value = 45
$ws.Cells.Item(194, 1).Value = '12:02:21'
$ws.Cells.Item(194, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(194, 4).Value = 32
$ws.Cells.Item(198, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(199, 3).Value = '17_179 Y 38'
$ws.Cells.Item(200, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(219, 3).Value = '215D_EL PATO'
$ws.Cells.Item(220, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(224, 1).Value = '13:18:32'
$ws.Cells.Item(224, 4).Value = 3
$ws.Cells.Item(226, 1).Value = '13:18:32'
$ws.Cells.Item(226, 4).Value = 8
$ws.Cells.Item(227, 1).Value = '13:18:32'
$ws.Cells.Item(227, 4).Value = 8
$ws.Cells.Item(229, 1).Value = '13:18:32'
$ws.Cells.Item(229, 2).Value = '13:33'
$ws.Cells.Item(229, 3).Value = '10_OLMOS'
$ws.Cells.Item(229, 4).Value = 15
$ws.Cells.Item(230, 1).Value = '13:18:32'
$ws.Cells.Item(230, 2).Value = '13:34'
$ws.Cells.Item(230, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(230, 4).Value = 16
$ws.Cells.Item(231, 1).Value = '13:18:32'
$ws.Cells.Item(231, 2).Value = '13:36'
$ws.Cells.Item(231, 3).Value = '15_ABASTO'
$ws.Cells.Item(231, 4).Value = 18
$ws.Cells.Item(232, 1).Value = '13:18:32'
$ws.Cells.Item(232, 2).Value = '13:46'
$ws.Cells.Item(232, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(232, 4).Value = 28
$ws.Cells.Item(233, 1).Value = '13:18:32'
$ws.Cells.Item(233, 2).Value = '13:46'
$ws.Cells.Item(233, 3).Value = '17_ROMERO'
$ws.Cells.Item(233, 4).Value = 28
$ws.Cells.Item(234, 1).Value = '13:18:32'
$ws.Cells.Item(234, 2).Value = '13:50'
$ws.Cells.Item(234, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(234, 4).Value = 32
$ws.Cells.Item(235, 1).Value = '13:18:32'
$ws.Cells.Item(235, 2).Value = '13:50'
$ws.Cells.Item(235, 3).Value = '215A_EL PATO'
$ws.Cells.Item(235, 4).Value = 32
$ws.Cells.Item(236, 1).Value = '12:37:00'
$ws.Cells.Item(236, 2).Value = '13:51'
$ws.Cells.Item(236, 3).Value = '215A_EL PATO'
$ws.Cells.Item(236, 4).Value = 74
$ws.Cells.Item(237, 1).Value = '13:18:32'
$ws.Cells.Item(237, 2).Value = '13:55'
$ws.Cells.Item(237, 3).Value = '16_SANTA ANA'
$ws.Cells.Item(237, 4).Value = 37
$ws.Cells.Item(238, 1).Value = '13:18:32'
$ws.Cells.Item(238, 2).Value = '13:56'
$ws.Cells.Item(238, 3).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(238, 4).Value = 38
$ws.Cells.Item(239, 1).Value = '13:18:32'
$ws.Cells.Item(239, 2).Value = '13:56'
$ws.Cells.Item(239, 3).Value = '225_GOMEZ'
$ws.Cells.Item(239, 4).Value = 38
$ws.Cells.Item(240, 1).Value = '12:37:00'
$ws.Cells.Item(240, 2).Value = '13:57'
$ws.Cells.Item(240, 3).Value = '16_P MOR-167 Y 521'
$ws.Cells.Item(240, 4).Value = 80
$ws.Cells.Item(241, 1).Value = '13:18:32'
$ws.Cells.Item(241, 2).Value = '14:04'
$ws.Cells.Item(241, 3).Value = '17_ROMERO'
$ws.Cells.Item(241, 4).Value = 46
$ws.Cells.Item(242, 1).Value = '13:18:32'
$ws.Cells.Item(242, 2).Value = '14:04'
$ws.Cells.Item(242, 3).Value = '23_HERNANDEZ'
$ws.Cells.Item(242, 4).Value = 46
$ws.Cells.Item(243, 1).Value = '13:18:32'
$ws.Cells.Item(243, 2).Value = '14:16'
$ws.Cells.Item(243, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(243, 4).Value = 58
$ws.Cells.Item(243, 5).Value = 'LP1912'
$ws.Cells.Item(244, 1).Value = '12:54:24'
$ws.Cells.Item(244, 2).Value = '14:17'
$ws.Cells.Item(244, 3).Value = '27_EL RETIRO'
$ws.Cells.Item(244, 4).Value = 83
$ws.Cells.Item(244, 5).Value = 'LP1912'
$ws.Cells.Item(245, 1).Value = '13:18:32'
$ws.Cells.Item(245, 2).Value = '14:20'
$ws.Cells.Item(245, 3).Value = '215C_EL PATO'
$ws.Cells.Item(245, 4).Value = 62
$ws.Cells.Item(245, 5).Value = 'LP1912'
$ws.Cells.Item(246, 1).Value = '13:18:32'
$ws.Cells.Item(246, 2).Value = '14:21'
$ws.Cells.Item(246, 3).Value = '26_HERNANDEZ'
$ws.Cells.Item(246, 4).Value = 63
$ws.Cells.Item(246, 5).Value = 'LP1912'
$ws.Cells.Item(247, 1).Value = '12:54:24'
$ws.Cells.Item(247, 2).Value = '14:39'
$ws.Cells.Item(247, 3).Value = '14_ABASTO'
$ws.Cells.Item(247, 4).Value = 105
$ws.Cells.Item(247, 5).Value = 'LP1912'
$ws.Cells.Item(248, 1).Value = '13:18:32'
$ws.Cells.Item(248, 2).Value = '14:44'
$ws.Cells.Item(248, 3).Value = '14_ABASTO'
$ws.Cells.Item(248, 4).Value = 86
$ws.Cells.Item(248, 5).Value = 'LP1912'
$ws.Cells.Item(249, 1).Value = '13:18:32'
$ws.Cells.Item(249, 2).Value = '14:56'
$ws.Cells.Item(249, 3).Value = '16_P MOR-SANTA ANA'
$ws.Cells.Item(249, 4).Value = 98
$ws.Cells.Item(249, 5).Value = 'LP1912'
$ws.Cells.Item(250, 1).Value = '13:18:32'
$ws.Cells.Item(250, 2).Value = '14:58'
$ws.Cells.Item(250, 3).Value = '215B_EL PATO'
$ws.Cells.Item(250, 4).Value = 100
$ws.Cells.Item(250, 5).Value = 'LP1912'
$ws.Cells.Item(251, 1).Value = '13:18:32'
$ws.Cells.Item(251, 2).Value = '15:00'
$ws.Cells.Item(251, 3).Value = '81_EL PELIGRO'
$ws.Cells.Item(251, 4).Value = 102
$ws.Cells.Item(251, 5).Value = 'LP1912'
$ws.Cells.Item(252, 1).Value = '13:18:32'
$ws.Cells.Item(252, 2).Value = '15:05'
$ws.Cells.Item(252, 3).Value = '10_OLMOS'
$ws.Cells.Item(252, 4).Value = 107
$ws.Cells.Item(252, 5).Value = 'LP1912'
$ws.Cells.Item(253, 1).Value = '13:18:32'
$ws.Cells.Item(253, 2).Value = '15:13'
$ws.Cells.Item(253, 3).Value = '11_ETCHEVERRY'
$ws.Cells.Item(253, 4).Value = 115
$ws.Cells.Item(253, 5).Value = 'LP1912'

# ---- Sheet: LP1912-215 ----
$ws = $wb.Worksheets.Item('LP1912-215')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 13:18:32'
$ws.Cells.Item(3, 1).Value = 'Total filas: 31'
$ws.Cells.Item(33, 1).Value = '13:18:32'
$ws.Cells.Item(33, 4).Value = 32
$ws.Cells.Item(35, 1).Value = '13:18:32'
$ws.Cells.Item(35, 4).Value = 62
$ws.Cells.Item(36, 1).Value = '13:18:32'
$ws.Cells.Item(36, 2).Value = '14:58'
$ws.Cells.Item(36, 3).Value = '215B_EL PATO'
$ws.Cells.Item(36, 4).Value = 100
$ws.Cells.Item(36, 5).Value = 'LP1912'

# ---- Sheet: 6203-6173 ----
$ws = $wb.Worksheets.Item('6203-6173')
$ws.Cells.Item(2, 1).Value = 'Última actualización: 13:18:32'
$ws.Cells.Item(38, 1).Value = '13:18:32'
$ws.Cells.Item(38, 4).Value = 13
$ws.Cells.Item(39, 1).Value = '13:18:32'
$ws.Cells.Item(39, 4).Value = 51
$ws.Cells.Item(40, 1).Value = '13:18:32'
$ws.Cells.Item(40, 4).Value = 95
